# Agrega la columna "Comuna" (rutas optimas por comuna) a la hoja Pedidos.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insertar una nueva columna en la posicion 17 (Q), desplazando
#    Motivo..Foto Enviado de Q..V a R..W.
$ws.Columns.Item(17).EntireColumn.Insert()

# 2) Encabezado de la nueva columna.
$ws.Range("Q1").Value2 = "Comuna"

# 3) Ajustar anchos de columna (en unidades OOXML "width", ColumnWidth = width - 0.83).
$ws.Columns.Item(17).ColumnWidth = 18 - 0.83   # Q - Comuna
$ws.Columns.Item(18).ColumnWidth = 15 - 0.83   # R - Motivo
$ws.Columns.Item(19).ColumnWidth = 18 - 0.83   # S - Estado
$ws.Columns.Item(20).ColumnWidth = 15 - 0.83   # T - Estado Pago
$ws.Columns.Item(21).ColumnWidth = 15 - 0.83   # U - Tipo Pedido
$ws.Columns.Item(22).ColumnWidth = 20 - 0.83   # V - Cobranza
$ws.Columns.Item(23).ColumnWidth = 25 - 0.83   # W - Foto Enviado

# 4) Valores de Comuna por pedido (derivados de la direccion de despacho).
$comunas = @{
    2  = "Las Condes"
    3  = "Providencia"
    4  = "Vitacura"
    5  = "Ñuñoa"
    6  = "La Reina"
    7  = "San Miguel"
    8  = "Maipú"
    9  = "Las Condes"
    10 = "Peñalolén"
}
foreach ($row in $comunas.Keys) {
    $ws.Range("Q$row").Value2 = $comunas[$row]
}

# 5) Precios de envio actualizados (columna L - Precio Envio) segun la nueva
#    tabla de precios por comuna.
$precios = @{
    2  = 7000
    3  = 10000
    4  = 7000
    5  = 15000
    6  = 15000
    7  = 25000
    8  = 30000
    9  = 7000
    10 = 25000
}
foreach ($row in $precios.Keys) {
    $ws.Range("L$row").Value2 = $precios[$row]
}
